$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0

$ws.Range("F3").Value = [double]"2.689122860863488e-09"
$ws.Range("G3").Value = 28.17704182552777
$ws.Range("H3").Value = 20.65093525419323

$ws.Range("G4").Value = 23.21271746461933
$ws.Range("I4").Value = 0

$ws.Range("G5").Value = 17.95850381196586
$ws.Range("H5").Value = 13.5511734331143

$ws.Range("C6").Value = [double]"2.689122860863488e-09"
$ws.Range("G6").Value = 18.54591470160023
$ws.Range("I6").Value = 8.255532014288399

$ws.Range("C7").Value = 28.17704182552777
$ws.Range("D7").Value = 23.21271746461933
$ws.Range("E7").Value = 17.95850381196586
$ws.Range("F7").Value = 18.54591470160023
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 20.65093525419323
$ws.Range("E8").Value = 13.5511734331143
$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 21.62916983090976

$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 8.255532014288399
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 21.62916983090976
